$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunTIMESReportScript")

# C15: drop the "What it was before: " prefix that used to live in C17,
# replacing the old single-run path with the per-case-name path.
$ws.Range("C15").Value = 'C:\veda\GAMS_WrkTIMES\%$case_name%\%$case_name%.lst'

# C17: now holds the new "max 1 runs" label.
$ws.Range("C17").Value = "If it is run with max 1 runs:"

# C18: new row added, holding the original single-run path text.
$ws.Range("C18").Value = 'C:\veda\GAMS_WrkTIMES\%$case_name%.lst'

# Scroll/selection state moved down to the newly added row.
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 10
